# Key Personnel List.xlsx — clean-up edit
#
# The "Full Title" / "Address" sub-header row was reworded:
#   - the cell that used to read "Full Title" now reads "LONG"
#   - the cell that used to read "Address"    now reads "ADDRESS"
#
# Both cells are the left-most cell of a merged A:G band (rows 2 and 3),
# so setting .Value on the single anchor cell is sufficient — Excel
# automatically reflects it across the merged range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "LONG"
$ws.Range("A3").Value = "ADDRESS"
